$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.70913233333333
$ws.Range("H2").Value = 32.127397
$ws.Range("I2").Value = 0.007451729107954897
$ws.Range("J2").Value = 0.007451729107954897
$ws.Range("M2").Value = 6.111751666666666
$ws.Range("N2").Value = 18.335255
$ws.Range("O2").Value = 0.6061514841909396
$ws.Range("P2").Value = 0.6061514841909394
$ws.Range("Q2").Value = 65.45155738680388
$ws.Range("R2").Value = 589.0640164812349
$ws.Range("S2").Value = 0.004516876658575687
$ws.Range("T2").Value = 0.004516876658575686
$ws.Range("G3").Value = 10.70913233333333
$ws.Range("H3").Value = 32.127397
$ws.Range("I3").Value = 0.007451729107954897
$ws.Range("J3").Value = 0.007451729107954897
$ws.Range("O3").Value = 0.2731664420559804
$ws.Range("P3").Value = 0.2731664420559804
$ws.Range("Q3").Value = 29.49620602222933
$ws.Range("R3").Value = 265.465854200064
$ws.Range("S3").Value = 0.002035562327585024
$ws.Range("T3").Value = 0.002035562327585024
$ws.Range("G4").Value = 10.70913233333333
$ws.Range("H4").Value = 32.127397
$ws.Range("I4").Value = 0.007451729107954897
$ws.Range("J4").Value = 0.007451729107954897
$ws.Range("M4").Value = 0.568439
$ws.Range("N4").Value = 1.705317
$ws.Range("O4").Value = 0.0563766596410053
$ws.Range("P4").Value = 0.05637665964100529
$ws.Range("Q4").Value = 6.087488474427667
$ws.Range("R4").Value = 54.787396269849
$ws.Range("S4").Value = 0.0004201035956561452
$ws.Range("T4").Value = 0.0004201035956561452
$ws.Range("G5").Value = 10.70913233333333
$ws.Range("H5").Value = 32.127397
$ws.Range("I5").Value = 0.007451729107954897
$ws.Range("J5").Value = 0.007451729107954897
$ws.Range("M5").Value = 0.3689163333333333
$ws.Range("N5").Value = 1.106749
$ws.Range("O5").Value = 0.036588394815171
$ws.Range("P5").Value = 0.036588394815171
$ws.Range("Q5").Value = 3.950773833594778
$ws.Range("R5").Value = 35.556964502353
$ws.Range("S5").Value = 0.0002726468066575558
$ws.Range("T5").Value = 0.0002726468066575558
$ws.Range("G6").Value = 10.70913233333333
$ws.Range("H6").Value = 32.127397
$ws.Range("I6").Value = 0.007451729107954897
$ws.Range("J6").Value = 0.007451729107954897
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2794673333333333
$ws.Range("N6").Value = 0.838402
$ws.Range("O6").Value = 0.02771701929690382
$ws.Range("P6").Value = 0.02771701929690381
$ws.Range("Q6").Value = 2.992852655510445
$ws.Range("R6").Value = 26.935673899594
$ws.Range("S6").Value = 0.0002065397194804857
$ws.Range("T6").Value = 0.0002065397194804857
$ws.Range("I7").Value = 0.03290895798513831
$ws.Range("J7").Value = 0.03290895798513832
$ws.Range("M7").Value = 6.111751666666666
$ws.Range("N7").Value = 18.335255
$ws.Range("O7").Value = 0.6061514841909396
$ws.Range("P7").Value = 0.6061514841909394
$ws.Range("Q7").Value = 289.0527179530472
$ws.Range("R7").Value = 2601.474461577424
$ws.Range("S7").Value = 0.01994781372586886
$ws.Range("T7").Value = 0.01994781372586886
$ws.Range("I8").Value = 0.03290895798513831
$ws.Range("J8").Value = 0.03290895798513832
$ws.Range("O8").Value = 0.2731664420559804
$ws.Range("P8").Value = 0.2731664420559804
$ws.Range("S8").Value = 0.008989622964569979
$ws.Range("T8").Value = 0.00898962296456998
$ws.Range("I9").Value = 0.03290895798513831
$ws.Range("J9").Value = 0.03290895798513832
$ws.Range("M9").Value = 0.568439
$ws.Range("N9").Value = 1.705317
$ws.Range("O9").Value = 0.0563766596410053
$ws.Range("P9").Value = 0.05637665964100529
$ws.Range("Q9").Value = 26.88408281322167
$ws.Range("R9").Value = 241.956745318995
$ws.Range("S9").Value = 0.001855297123468286
$ws.Range("T9").Value = 0.001855297123468286
$ws.Range("I10").Value = 0.03290895798513831
$ws.Range("J10").Value = 0.03290895798513832
$ws.Range("M10").Value = 0.3689163333333333
$ws.Range("N10").Value = 1.106749
$ws.Range("O10").Value = 0.036588394815171
$ws.Range("P10").Value = 0.036588394815171
$ws.Range("Q10").Value = 17.44774242527944
$ws.Range("R10").Value = 157.029681827515
$ws.Range("S10").Value = 0.001204085947716115
$ws.Range("T10").Value = 0.001204085947716115
$ws.Range("I11").Value = 0.03290895798513831
$ws.Range("J11").Value = 0.03290895798513832
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2794673333333333
$ws.Range("N11").Value = 0.838402
$ws.Range("O11").Value = 0.02771701929690382
$ws.Range("P11").Value = 0.02771701929690381
$ws.Range("Q11").Value = 13.21728968794111
$ws.Range("R11").Value = 118.95560719147
$ws.Range("S11").Value = 0.0009121382235150755
$ws.Range("T11").Value = 0.0009121382235150756
$ws.Range("G12").Value = 411.37678
$ws.Range("H12").Value = 1234.13034
$ws.Range("I12").Value = 0.2862480573072345
$ws.Range("J12").Value = 0.2862480573072345
$ws.Range("M12").Value = 6.111751666666666
$ws.Range("N12").Value = 18.335255
$ws.Range("O12").Value = 0.6061514841909396
$ws.Range("P12").Value = 0.6061514841909394
$ws.Range("Q12").Value = 2514.232720792966
$ws.Range("R12").Value = 22628.0944871367
$ws.Range("S12").Value = 0.1735096847835533
$ws.Range("T12").Value = 0.1735096847835532
$ws.Range("G13").Value = 411.37678
$ws.Range("H13").Value = 1234.13034
$ws.Range("I13").Value = 0.2862480573072345
$ws.Range("J13").Value = 0.2862480573072345
$ws.Range("O13").Value = 0.2731664420559804
$ws.Range("P13").Value = 0.2731664420559804
$ws.Range("Q13").Value = 1133.05671066112
$ws.Range("R13").Value = 10197.51039595008
$ws.Range("S13").Value = 0.07819336336005361
$ws.Range("T13").Value = 0.07819336336005361
$ws.Range("G14").Value = 411.37678
$ws.Range("H14").Value = 1234.13034
$ws.Range("I14").Value = 0.2862480573072345
$ws.Range("J14").Value = 0.2862480573072345
$ws.Range("M14").Value = 0.568439
$ws.Range("N14").Value = 1.705317
$ws.Range("O14").Value = 0.0563766596410053
$ws.Range("P14").Value = 0.05637665964100529
$ws.Range("Q14").Value = 233.84260544642
$ws.Range("R14").Value = 2104.58344901778
$ws.Range("S14").Value = 0.01613770929970894
$ws.Range("T14").Value = 0.01613770929970893
$ws.Range("G15").Value = 411.37678
$ws.Range("H15").Value = 1234.13034
$ws.Range("I15").Value = 0.2862480573072345
$ws.Range("J15").Value = 0.2862480573072345
$ws.Range("M15").Value = 0.3689163333333333
$ws.Range("N15").Value = 1.106749
$ws.Range("O15").Value = 0.036588394815171
$ws.Range("P15").Value = 0.036588394815171
$ws.Range("Q15").Value = 151.7636132960733
$ws.Range("R15").Value = 1365.87251966466
$ws.Range("S15").Value = 0.01047335693583279
$ws.Range("T15").Value = 0.01047335693583279
$ws.Range("G16").Value = 411.37678
$ws.Range("H16").Value = 1234.13034
$ws.Range("I16").Value = 0.2862480573072345
$ws.Range("J16").Value = 0.2862480573072345
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2794673333333333
$ws.Range("N16").Value = 0.838402
$ws.Range("O16").Value = 0.02771701929690382
$ws.Range("P16").Value = 0.02771701929690381
$ws.Range("Q16").Value = 114.9663717018533
$ws.Range("R16").Value = 1034.69734531668
$ws.Range("S16").Value = 0.007933942928085847
$ws.Range("T16").Value = 0.007933942928085845
$ws.Range("G17").Value = 173.2560603333334
$ws.Range("H17").Value = 519.768181
$ws.Range("I17").Value = 0.12055666021578
$ws.Range("J17").Value = 0.12055666021578
$ws.Range("M17").Value = 6.111751666666666
$ws.Range("N17").Value = 18.335255
$ws.Range("O17").Value = 0.6061514841909396
$ws.Range("P17").Value = 0.6061514841909394
$ws.Range("Q17").Value = 1058.898015502351
$ws.Range("R17").Value = 9530.082139521153
$ws.Range("S17").Value = 0.07307559851889783
$ws.Range("T17").Value = 0.07307559851889782
$ws.Range("G18").Value = 173.2560603333334
$ws.Range("H18").Value = 519.768181
$ws.Range("I18").Value = 0.12055666021578
$ws.Range("J18").Value = 0.12055666021578
$ws.Range("O18").Value = 0.2731664420559804
$ws.Range("P18").Value = 0.2731664420559804
$ws.Range("Q18").Value = 477.1998600003413
$ws.Range("R18").Value = 4294.798740003072
$ws.Range("S18").Value = 0.03293203393729638
$ws.Range("T18").Value = 0.03293203393729638
$ws.Range("G19").Value = 173.2560603333334
$ws.Range("H19").Value = 519.768181
$ws.Range("I19").Value = 0.12055666021578
$ws.Range("J19").Value = 0.12055666021578
$ws.Range("M19").Value = 0.568439
$ws.Range("N19").Value = 1.705317
$ws.Range("O19").Value = 0.0563766596410053
$ws.Range("P19").Value = 0.05637665964100529
$ws.Range("Q19").Value = 98.48550167981968
$ws.Range("R19").Value = 886.3695151183771
$ws.Range("S19").Value = 0.006796581800441352
$ws.Range("T19").Value = 0.006796581800441352
$ws.Range("G20").Value = 173.2560603333334
$ws.Range("H20").Value = 519.768181
$ws.Range("I20").Value = 0.12055666021578
$ws.Range("J20").Value = 0.12055666021578
$ws.Range("M20").Value = 0.3689163333333333
$ws.Range("N20").Value = 1.106749
$ws.Range("O20").Value = 0.036588394815171
$ws.Range("P20").Value = 0.036588394815171
$ws.Range("Q20").Value = 63.91699050595212
$ws.Range("R20").Value = 575.252914553569
$ws.Range("S20").Value = 0.004410974681573377
$ws.Range("T20").Value = 0.004410974681573377
$ws.Range("G21").Value = 173.2560603333334
$ws.Range("H21").Value = 519.768181
$ws.Range("I21").Value = 0.12055666021578
$ws.Range("J21").Value = 0.12055666021578
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.2794673333333333
$ws.Range("N21").Value = 0.838402
$ws.Range("O21").Value = 0.02771701929690382
$ws.Range("P21").Value = 0.02771701929690381
$ws.Range("Q21").Value = 48.41940916519579
$ws.Range("R21").Value = 435.774682486762
$ws.Range("S21").Value = 0.00334147127757105
$ws.Range("T21").Value = 0.00334147127757105
$ws.Range("G22").Value = 794.4973246666667
$ws.Range("H22").Value = 2383.491974
$ws.Range("I22").Value = 0.5528345953838922
$ws.Range("J22").Value = 0.5528345953838923
$ws.Range("M22").Value = 6.111751666666666
$ws.Range("N22").Value = 18.335255
$ws.Range("O22").Value = 0.6061514841909396
$ws.Range("P22").Value = 0.6061514841909394
$ws.Range("Q22").Value = 4855.770348193708
$ws.Range("R22").Value = 43701.93313374336
$ws.Range("S22").Value = 0.3351015105040438
$ws.Range("T22").Value = 0.3351015105040439
$ws.Range("G23").Value = 794.4973246666667
$ws.Range("H23").Value = 2383.491974
$ws.Range("I23").Value = 0.5528345953838922
$ws.Range("J23").Value = 0.5528345953838923
$ws.Range("O23").Value = 0.2731664420559804
$ws.Range("P23").Value = 0.2731664420559804
$ws.Range("Q23").Value = 2188.287159318699
$ws.Range("R23").Value = 19694.58443386829
$ws.Range("S23").Value = 0.1510158594664754
$ws.Range("T23").Value = 0.1510158594664754
$ws.Range("G24").Value = 794.4973246666667
$ws.Range("H24").Value = 2383.491974
$ws.Range("I24").Value = 0.5528345953838922
$ws.Range("J24").Value = 0.5528345953838923
$ws.Range("M24").Value = 0.568439
$ws.Range("N24").Value = 1.705317
$ws.Range("O24").Value = 0.0563766596410053
$ws.Range("P24").Value = 0.05637665964100529
$ws.Range("Q24").Value = 451.6232647361954
$ws.Range("R24").Value = 4064.609382625758
$ws.Range("S24").Value = 0.03116696782173057
$ws.Range("T24").Value = 0.03116696782173058
$ws.Range("G25").Value = 794.4973246666667
$ws.Range("H25").Value = 2383.491974
$ws.Range("I25").Value = 0.5528345953838922
$ws.Range("J25").Value = 0.5528345953838923
$ws.Range("M25").Value = 0.3689163333333333
$ws.Range("N25").Value = 1.106749
$ws.Range("O25").Value = 0.036588394815171
$ws.Range("P25").Value = 0.036588394815171
$ws.Range("Q25").Value = 293.1030398591696
$ws.Range("R25").Value = 2637.927358732526
$ws.Range("S25").Value = 0.02022733044339116
$ws.Range("T25").Value = 0.02022733044339117
$ws.Range("G26").Value = 794.4973246666667
$ws.Range("H26").Value = 2383.491974
$ws.Range("I26").Value = 0.5528345953838922
$ws.Range("J26").Value = 0.5528345953838923
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.2794673333333333
$ws.Range("N26").Value = 0.838402
$ws.Range("O26").Value = 0.02771701929690382
$ws.Range("P26").Value = 0.02771701929690381
$ws.Range("Q26").Value = 222.0360486650609
$ws.Range("R26").Value = 1998.324437985548
$ws.Range("S26").Value = 0.01532292714825135
$ws.Range("T26").Value = 0.01532292714825136
